$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 514.4
$ws.Range("I12").Value = 670.2857
$ws.Range("J12").Value = 150.66667
$ws.Range("K12").Value = 670.2857
$ws.Range("L12").Value = 150.66667
$ws.Range("M12").Value = -500.2857
$ws.Range("N12").Value = -490.66667
$ws.Range("H96").Value = 91899.27
$ws.Range("I96").Value = 167069.67
$ws.Range("J96").Value = 1694.8
$ws.Range("K96").Value = 501209.01
$ws.Range("L96").Value = 5084.4
$ws.Range("M96").Value = -499836.01
$ws.Range("N96").Value = -7830.4
$ws.Range("H103").Value = 1051
$ws.Range("H132").Value = 4661.548
$ws.Range("I132").Value = 4838.3613
$ws.Range("J132").Value = 3600.6667
$ws.Range("K132").Value = 14515.0839
$ws.Range("L132").Value = 10802.0001
$ws.Range("M132").Value = -11985.0839
$ws.Range("N132").Value = -15862.0001
$ws.Range("H135").Value = 1377.1428
$ws.Range("I135").Value = 934.5454999999999
$ws.Range("K135").Value = 8410.9095
$ws.Range("M135").Value = -5875.9095
$ws.Range("H137").Value = 6859.3335
$ws.Range("J137").Value = 7237
$ws.Range("L137").Value = 21711
$ws.Range("N137").Value = -26811
$ws.Range("H138").Value = 7056.7607
$ws.Range("I138").Value = 6782.1113
$ws.Range("J138").Value = 7123.5674
$ws.Range("K138").Value = 20346.3339
$ws.Range("L138").Value = 21370.7022
$ws.Range("M138").Value = -15206.3339
$ws.Range("N138").Value = -31650.7022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 254194.31
$ws.Range("J45").Value = 1578.25
$ws.Range("L45").Value = 1578.25
$ws.Range("N45").Value = -2332.25
$ws.Range("H61").Value = 627413.25
$ws.Range("I61").Value = 2615.1428
$ws.Range("J61").Value = 5001000
$ws.Range("K61").Value = 2615.1428
$ws.Range("L61").Value = 5001000
$ws.Range("M61").Value = -2403.1428
$ws.Range("N61").Value = -5001424
$ws.Range("H74").Value = 2654.5642
$ws.Range("I74").Value = 2044.88
$ws.Range("J74").Value = 3743.2856
$ws.Range("K74").Value = 2044.88
$ws.Range("L74").Value = 3743.2856
$ws.Range("M74").Value = -1170.88
$ws.Range("N74").Value = -5491.2856
$ws.Range("H77").Value = 2654.5642
$ws.Range("I77").Value = 2044.88
$ws.Range("J77").Value = 3743.2856
$ws.Range("K77").Value = 10224.4
$ws.Range("L77").Value = 18716.428
$ws.Range("M77").Value = -5856.400000000001
$ws.Range("N77").Value = -27452.428
$ws.Range("H102").Value = 5935447
$ws.Range("J102").Value = 126624.5
$ws.Range("L102").Value = 126624.5
$ws.Range("N102").Value = -129868.5
$ws.Range("H110").Value = 34484744
$ws.Range("J110").Value = 3610.25
$ws.Range("L110").Value = 3610.25
$ws.Range("N110").Value = -7700.25
$ws.Range("H122").Value = 2575.0908
$ws.Range("I122").Value = 1371.0834
$ws.Range("J122").Value = 4019.9
$ws.Range("K122").Value = 4113.2502
$ws.Range("L122").Value = 12059.7
$ws.Range("M122").Value = -1663.2502
$ws.Range("N122").Value = -16959.7
$ws.Range("H132").Value = 17934228
$ws.Range("I132").Value = 19617946
$ws.Range("J132").Value = 3341999.2
$ws.Range("K132").Value = 58853838
$ws.Range("L132").Value = 10025997.6
$ws.Range("M132").Value = -58851308
$ws.Range("N132").Value = -10031057.6
$ws.Range("H136").Value = 627413.25
$ws.Range("I136").Value = 2615.1428
$ws.Range("J136").Value = 5001000
$ws.Range("K136").Value = 7845.428400000001
$ws.Range("L136").Value = 15003000
$ws.Range("M136").Value = -5295.428400000001
$ws.Range("N136").Value = -15008100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4466.952
$ws.Range("I134").Value = 4115.2856
$ws.Range("J134").Value = 5170.2856
$ws.Range("K134").Value = 12345.8568
$ws.Range("L134").Value = 15510.8568
$ws.Range("M134").Value = -9810.856800000001
$ws.Range("N134").Value = -20580.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4072.1292
$ws.Range("I31").Value = 2717.6365
$ws.Range("J31").Value = 4817.1
$ws.Range("K31").Value = 2717.6365
$ws.Range("L31").Value = 4817.1
$ws.Range("M31").Value = -2422.6365
$ws.Range("N31").Value = -5407.1
$ws.Range("H34").Value = 4072.1292
$ws.Range("I34").Value = 2717.6365
$ws.Range("J34").Value = 4817.1
$ws.Range("K34").Value = 2717.6365
$ws.Range("L34").Value = 4817.1
$ws.Range("M34").Value = -2515.6365
$ws.Range("N34").Value = -5221.1
$ws.Range("H58").Value = 215977.66
$ws.Range("I58").Value = 1877.55
$ws.Range("J58").Value = 374570.34
$ws.Range("K58").Value = 1877.55
$ws.Range("L58").Value = 374570.34
$ws.Range("M58").Value = -1674.55
$ws.Range("N58").Value = -374976.34
$ws.Range("H122").Value = 1768.6666
$ws.Range("I122").Value = 1322.8
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 3968.4
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -1518.4
$ws.Range("N122").Value = -16894
$ws.Range("H132").Value = 670209.6
$ws.Range("I132").Value = 3742
$ws.Range("K132").Value = 11226
$ws.Range("M132").Value = -8696
$ws.Range("H136").Value = 215977.66
$ws.Range("I136").Value = 1877.55
$ws.Range("J136").Value = 374570.34
$ws.Range("K136").Value = 5632.65
$ws.Range("L136").Value = 1123711.02
$ws.Range("M136").Value = -3082.65
$ws.Range("N136").Value = -1128811.02

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 754.1429000000001
$ws.Range("J5").Value = 950
$ws.Range("L5").Value = 2850
$ws.Range("N5").Value = -3074
$ws.Range("H56").Value = 5199.25
$ws.Range("I56").Value = 5199.25
$ws.Range("K56").Value = 5199.25
$ws.Range("M56").Value = -4669.25
$ws.Range("H81").Value = 10256.429
$ws.Range("I81").Value = 2149.5
$ws.Range("J81").Value = 13499.2
$ws.Range("K81").Value = 6448.5
$ws.Range("L81").Value = 40497.60000000001
$ws.Range("M81").Value = -5325.5
$ws.Range("N81").Value = -42743.60000000001
$ws.Range("H84").Value = 10256.429
$ws.Range("I84").Value = 2149.5
$ws.Range("J84").Value = 13499.2
$ws.Range("K84").Value = 19345.5
$ws.Range("L84").Value = 121492.8
$ws.Range("M84").Value = -13729.5
$ws.Range("N84").Value = -132724.8
$ws.Range("H113").Value = 2333265.2
$ws.Range("J113").Value = 2851324.2
$ws.Range("L113").Value = 8553972.600000001
$ws.Range("N113").Value = -8558312.600000001
$ws.Range("H132").Value = 2348.3076
$ws.Range("J132").Value = 2583.2856
$ws.Range("L132").Value = 23249.5704
$ws.Range("N132").Value = -28309.5704
$ws.Range("H135").Value = 754.1429000000001
$ws.Range("J135").Value = 950
$ws.Range("L135").Value = 8550
$ws.Range("N135").Value = -13620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125002250
$ws.Range("I80").Value = 2334.5
$ws.Range("J80").Value = 500002000
$ws.Range("K80").Value = 2334.5
$ws.Range("L80").Value = 500002000
$ws.Range("M80").Value = -1336.5
$ws.Range("N80").Value = -500003996
$ws.Range("H83").Value = 125002250
$ws.Range("I83").Value = 2334.5
$ws.Range("J83").Value = 500002000
$ws.Range("K83").Value = 11672.5
$ws.Range("L83").Value = 2500010000
$ws.Range("M83").Value = -6680.5
$ws.Range("N83").Value = -2500019984
$ws.Range("H97").Value = 780.8946999999999
$ws.Range("I97").Value = 823.4375
$ws.Range("J97").Value = 554
$ws.Range("K97").Value = 823.4375
$ws.Range("L97").Value = 554
$ws.Range("M97").Value = -327.4375
$ws.Range("N97").Value = -1546
$ws.Range("H102").Value = 1984.8085
$ws.Range("I102").Value = 859.76666
$ws.Range("J102").Value = 3970.1765
$ws.Range("K102").Value = 859.76666
$ws.Range("L102").Value = 3970.1765
$ws.Range("M102").Value = 762.23334
$ws.Range("N102").Value = -7214.1765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 100003800
$ws.Range("I7").Value = 250002000
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 250002000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -250001888
$ws.Range("N7").Value = -5224
$ws.Range("H40").Value = 5191.706
$ws.Range("I40").Value = 4756.4
$ws.Range("J40").Value = 5373.0835
$ws.Range("K40").Value = 4756.4
$ws.Range("L40").Value = 5373.0835
$ws.Range("M40").Value = -4620.4
$ws.Range("N40").Value = -5645.0835
$ws.Range("H122").Value = 10916.667
$ws.Range("J122").Value = 6500
$ws.Range("L122").Value = 19500
$ws.Range("N122").Value = -24400
$ws.Range("H126").Value = 100003800
$ws.Range("I126").Value = 250002000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 750006000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -750003530
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 177488.33
$ws.Range("I132").Value = 289454.16
$ws.Range("J132").Value = 7105.522
$ws.Range("K132").Value = 868362.48
$ws.Range("L132").Value = 21316.566
$ws.Range("M132").Value = -865832.48
$ws.Range("N132").Value = -26376.566
$ws.Range("H136").Value = 25647162
$ws.Range("J136").Value = 6746.4
$ws.Range("L136").Value = 20239.2
$ws.Range("N136").Value = -25339.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 360.66666
$ws.Range("I107").Value = 381.30768
$ws.Range("J107").Value = 307
$ws.Range("K107").Value = 1143.92304
$ws.Range("L107").Value = 921
$ws.Range("M107").Value = 776.0769599999999
$ws.Range("N107").Value = -4761
$ws.Range("H126").Value = 8800.4
$ws.Range("I126").Value = 5002
$ws.Range("J126").Value = 11332.667
$ws.Range("K126").Value = 15006
$ws.Range("L126").Value = 33998.001
$ws.Range("M126").Value = -12536
$ws.Range("N126").Value = -38938.001
$ws.Range("H132").Value = 307921.5
$ws.Range("I132").Value = 327351.9
$ws.Range("K132").Value = 982055.7000000001
$ws.Range("M132").Value = -979525.7000000001
$ws.Range("H136").Value = 10691.36
$ws.Range("I136").Value = 10990.272
$ws.Range("J136").Value = 8499.333000000001
$ws.Range("K136").Value = 32970.81600000001
$ws.Range("L136").Value = 25497.999
$ws.Range("M136").Value = -30420.81600000001
$ws.Range("N136").Value = -30597.999
